$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 14:46"

# Refresh the full data block (A4:H199) with updated case counts
# and the re-sorted (descending by Casos totales) row order.
$data = New-Object 'object[,]' 196,8
$data[0,0] = "China"
$data[0,1] = 81171
$data[0,2] = 78
$data[0,3] = 73159
$data[0,4] = 4735
$data[0,5] = 1573
$data[0,6] = 7
$data[0,7] = 3277
$data[1,0] = "Italia"
$data[1,1] = 63927
$data[1,2] = 0
$data[1,3] = 7432
$data[1,4] = 50418
$data[1,5] = 3204
$data[1,6] = 0
$data[1,7] = 6077
$data[2,0] = "Estados Unidos"
$data[2,1] = 46168
$data[2,2] = 2434
$data[2,3] = 295
$data[2,4] = 45291
$data[2,5] = 1040
$data[2,6] = 29
$data[2,7] = 582
$data[3,0] = "España"
$data[3,1] = 39673
$data[3,2] = 4537
$data[3,3] = 3794
$data[3,4] = 33183
$data[3,5] = 2355
$data[3,6] = 385
$data[3,7] = 2696
$data[4,0] = "Alemania"
$data[4,1] = 31260
$data[4,2] = 2204
$data[4,3] = 749
$data[4,4] = 30379
$data[4,5] = 23
$data[4,6] = 9
$data[4,7] = 132
$data[5,0] = "Iran"
$data[5,1] = 24811
$data[5,2] = 1762
$data[5,3] = 8913
$data[5,4] = 13964
$data[5,5] = 0
$data[5,6] = 122
$data[5,7] = 1934
$data[6,0] = "Francia"
$data[6,1] = 19856
$data[6,2] = 0
$data[6,3] = 2200
$data[6,4] = 16796
$data[6,5] = 2082
$data[6,6] = 0
$data[6,7] = 860
$data[7,0] = "Suiza"
$data[7,1] = 9117
$data[7,2] = 322
$data[7,3] = 131
$data[7,4] = 8864
$data[7,5] = 141
$data[7,6] = 2
$data[7,7] = 122
$data[8,0] = "Corea del Sur"
$data[8,1] = 9037
$data[8,2] = 76
$data[8,3] = 3507
$data[8,4] = 5410
$data[8,5] = 59
$data[8,6] = 9
$data[8,7] = 120
$data[9,0] = "Reino Unido"
$data[9,1] = 6650
$data[9,2] = 0
$data[9,3] = 135
$data[9,4] = 6180
$data[9,5] = 20
$data[9,6] = 0
$data[9,7] = 335
$data[10,0] = "Paises Bajos"
$data[10,1] = 5560
$data[10,2] = 811
$data[10,3] = 2
$data[10,4] = 5282
$data[10,5] = 435
$data[10,6] = 63
$data[10,7] = 276
$data[11,0] = "Austria"
$data[11,1] = 4876
$data[11,2] = 402
$data[11,3] = 9
$data[11,4] = 4842
$data[11,5] = 19
$data[11,6] = 4
$data[11,7] = 25
$data[12,0] = "Belgica"
$data[12,1] = 4269
$data[12,2] = 526
$data[12,3] = 461
$data[12,4] = 3686
$data[12,5] = 381
$data[12,6] = 34
$data[12,7] = 122
$data[13,0] = "Noruega"
$data[13,1] = 2715
$data[13,2] = 90
$data[13,3] = 6
$data[13,4] = 2697
$data[13,5] = 44
$data[13,6] = 2
$data[13,7] = 12
$data[14,0] = "Portugal"
$data[14,1] = 2362
$data[14,2] = 302
$data[14,3] = 22
$data[14,4] = 2311
$data[14,5] = 47
$data[14,6] = 6
$data[14,7] = 29
$data[15,0] = "Suecia"
$data[15,1] = 2272
$data[15,2] = 226
$data[15,3] = 16
$data[15,4] = 2220
$data[15,5] = 115
$data[15,6] = 9
$data[15,7] = 36
$data[16,0] = "Australia"
$data[16,1] = 2144
$data[16,2] = 257
$data[16,3] = 118
$data[16,4] = 2018
$data[16,5] = 11
$data[16,6] = 1
$data[16,7] = 8
$data[17,0] = "Canada"
$data[17,1] = 2091
$data[17,2] = 0
$data[17,3] = 112
$data[17,4] = 1955
$data[17,5] = 1
$data[17,6] = 0
$data[17,7] = 24
$data[18,0] = "Brasil"
$data[18,1] = 1965
$data[18,2] = 41
$data[18,3] = 2
$data[18,4] = 1929
$data[18,5] = 18
$data[18,6] = 0
$data[18,7] = 34
$data[19,0] = "Israel"
$data[19,1] = 1656
$data[19,2] = 214
$data[19,3] = 49
$data[19,4] = 1605
$data[19,5] = 31
$data[19,6] = 1
$data[19,7] = 2
$data[20,0] = "Malasia"
$data[20,1] = 1624
$data[20,2] = 106
$data[20,3] = 183
$data[20,4] = 1426
$data[20,5] = 64
$data[20,6] = 1
$data[20,7] = 15
$data[21,0] = "Dinamarca"
$data[21,1] = 1577
$data[21,2] = 117
$data[21,3] = 1
$data[21,4] = 1544
$data[21,5] = 69
$data[21,6] = 8
$data[21,7] = 32
$data[22,0] = "Turquia"
$data[22,1] = 1529
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 1492
$data[22,5] = 0
$data[22,6] = 0
$data[22,7] = 37
$data[23,0] = "Chequia"
$data[23,1] = 1289
$data[23,2] = 53
$data[23,3] = 8
$data[23,4] = 1279
$data[23,5] = 19
$data[23,6] = 1
$data[23,7] = 2
$data[24,0] = "Japon"
$data[24,1] = 1140
$data[24,2] = 12
$data[24,3] = 285
$data[24,4] = 813
$data[24,5] = 54
$data[24,6] = 0
$data[24,7] = 42
$data[25,0] = "Irlanda"
$data[25,1] = 1125
$data[25,2] = 0
$data[25,3] = 5
$data[25,4] = 1114
$data[25,5] = 29
$data[25,6] = 0
$data[25,7] = 6
$data[26,0] = "Ecuador"
$data[26,1] = 981
$data[26,2] = 0
$data[26,3] = 3
$data[26,4] = 960
$data[26,5] = 2
$data[26,6] = 0
$data[26,7] = 18
$data[27,0] = "Chile"
$data[27,1] = 922
$data[27,2] = 176
$data[27,3] = 17
$data[27,4] = 903
$data[27,5] = 7
$data[27,6] = 0
$data[27,7] = 2
$data[28,0] = "Pakistan"
$data[28,1] = 918
$data[28,2] = 43
$data[28,3] = 13
$data[28,4] = 898
$data[28,5] = 0
$data[28,6] = 1
$data[28,7] = 7
$data[29,0] = "Luxemburgo"
$data[29,1] = 875
$data[29,2] = 0
$data[29,3] = 6
$data[29,4] = 861
$data[29,5] = 3
$data[29,6] = 0
$data[29,7] = 8
$data[30,0] = "Tailandia"
$data[30,1] = 827
$data[30,2] = 106
$data[30,3] = 52
$data[30,4] = 771
$data[30,5] = 7
$data[30,6] = 3
$data[30,7] = 4
$data[31,0] = "Polonia"
$data[31,1] = 799
$data[31,2] = 50
$data[31,3] = 1
$data[31,4] = 789
$data[31,5] = 3
$data[31,6] = 1
$data[31,7] = 9
$data[32,0] = "Finlandia"
$data[32,1] = 792
$data[32,2] = 92
$data[32,3] = 10
$data[32,4] = 781
$data[32,5] = 11
$data[32,6] = 0
$data[32,7] = 1
$data[33,0] = "Arabia Saudita"
$data[33,1] = 767
$data[33,2] = 205
$data[33,3] = 28
$data[33,4] = 738
$data[33,5] = 0
$data[33,6] = 1
$data[33,7] = 1
$data[34,0] = "Rumania"
$data[34,1] = 762
$data[34,2] = 186
$data[34,3] = 79
$data[34,4] = 675
$data[34,5] = 15
$data[34,6] = 1
$data[34,7] = 8
$data[35,0] = "Crucero"
$data[35,1] = 712
$data[35,2] = 0
$data[35,3] = 587
$data[35,4] = 115
$data[35,5] = 15
$data[35,6] = 2
$data[35,7] = 10
$data[36,0] = "Grecia"
$data[36,1] = 695
$data[36,2] = 0
$data[36,3] = 29
$data[36,4] = 647
$data[36,5] = 35
$data[36,6] = 2
$data[36,7] = 19
$data[37,0] = "Indonesia"
$data[37,1] = 686
$data[37,2] = 107
$data[37,3] = 30
$data[37,4] = 601
$data[37,5] = 0
$data[37,6] = 6
$data[37,7] = 55
$data[38,0] = "Islandia"
$data[38,1] = 648
$data[38,2] = 60
$data[38,3] = 51
$data[38,4] = 595
$data[38,5] = 13
$data[38,6] = 1
$data[38,7] = 2
$data[39,0] = "Singapur"
$data[39,1] = 558
$data[39,2] = 49
$data[39,3] = 156
$data[39,4] = 400
$data[39,5] = 14
$data[39,6] = 0
$data[39,7] = 2
$data[40,0] = "Sudafrica"
$data[40,1] = 554
$data[40,2] = 152
$data[40,3] = 4
$data[40,4] = 550
$data[40,5] = 2
$data[40,6] = 0
$data[40,7] = 0
$data[41,0] = "Filipinas"
$data[41,1] = 552
$data[41,2] = 90
$data[41,3] = 20
$data[41,4] = 497
$data[41,5] = 1
$data[41,6] = 2
$data[41,7] = 35
$data[42,0] = "India"
$data[42,1] = 519
$data[42,2] = 20
$data[42,3] = 40
$data[42,4] = 469
$data[42,5] = 0
$data[42,6] = 0
$data[42,7] = 10
$data[43,0] = "Catar"
$data[43,1] = 501
$data[43,2] = 0
$data[43,3] = 37
$data[43,4] = 464
$data[43,5] = 6
$data[43,6] = 0
$data[43,7] = 0
$data[44,0] = "Rusia"
$data[44,1] = 495
$data[44,2] = 57
$data[44,3] = 22
$data[44,4] = 472
$data[44,5] = 8
$data[44,6] = 0
$data[44,7] = 1
$data[45,0] = "Eslovenia"
$data[45,1] = 480
$data[45,2] = 38
$data[45,3] = 3
$data[45,4] = 473
$data[45,5] = 12
$data[45,6] = 1
$data[45,7] = 4
$data[46,0] = "Peru"
$data[46,1] = 395
$data[46,2] = 0
$data[46,3] = 1
$data[46,4] = 389
$data[46,5] = 19
$data[46,6] = 0
$data[46,7] = 5
$data[47,0] = "Barein"
$data[47,1] = 390
$data[47,2] = 13
$data[47,3] = 164
$data[47,4] = 223
$data[47,5] = 2
$data[47,6] = 1
$data[47,7] = 3
$data[48,0] = "Hong Kong"
$data[48,1] = 386
$data[48,2] = 29
$data[48,3] = 102
$data[48,4] = 280
$data[48,5] = 4
$data[48,6] = 0
$data[48,7] = 4
$data[49,0] = "Estonia"
$data[49,1] = 369
$data[49,2] = 17
$data[49,3] = 7
$data[49,4] = 362
$data[49,5] = 4
$data[49,6] = 0
$data[49,7] = 0
$data[50,0] = "Mexico"
$data[50,1] = 367
$data[50,2] = 51
$data[50,3] = 4
$data[50,4] = 359
$data[50,5] = 1
$data[50,6] = 1
$data[50,7] = 4
$data[51,0] = "Egipto"
$data[51,1] = 366
$data[51,2] = 0
$data[51,3] = 68
$data[51,4] = 279
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 19
$data[52,0] = "Croacia"
$data[52,1] = 361
$data[52,2] = 46
$data[52,3] = 5
$data[52,4] = 355
$data[52,5] = 6
$data[52,6] = 0
$data[52,7] = 1
$data[53,0] = "Panama"
$data[53,1] = 345
$data[53,2] = 0
$data[53,3] = 1
$data[53,4] = 338
$data[53,5] = 33
$data[53,6] = 0
$data[53,7] = 6
$data[54,0] = "Irak"
$data[54,1] = 316
$data[54,2] = 50
$data[54,3] = 75
$data[54,4] = 214
$data[54,5] = 0
$data[54,6] = 4
$data[54,7] = 27
$data[55,0] = "Colombia"
$data[55,1] = 306
$data[55,2] = 29
$data[55,3] = 6
$data[55,4] = 297
$data[55,5] = 0
$data[55,6] = 0
$data[55,7] = 3
$data[56,0] = "Libano"
$data[56,1] = 304
$data[56,2] = 37
$data[56,3] = 8
$data[56,4] = 292
$data[56,5] = 4
$data[56,6] = 0
$data[56,7] = 4
$data[57,0] = "Argentina"
$data[57,1] = 301
$data[57,2] = 0
$data[57,3] = 51
$data[57,4] = 246
$data[57,5] = 0
$data[57,6] = 0
$data[57,7] = 4
$data[58,0] = "Serbia"
$data[58,1] = 249
$data[58,2] = 0
$data[58,3] = 3
$data[58,4] = 243
$data[58,5] = 12
$data[58,6] = 0
$data[58,7] = 3
$data[59,0] = "Republica Dominicana"
$data[59,1] = 245
$data[59,2] = 0
$data[59,3] = 3
$data[59,4] = 239
$data[59,5] = 0
$data[59,6] = 0
$data[59,7] = 3
$data[60,0] = "Armenia"
$data[60,1] = 235
$data[60,2] = 0
$data[60,3] = 2
$data[60,4] = 233
$data[60,5] = 6
$data[60,6] = 0
$data[60,7] = 0
$data[61,0] = "Argelia"
$data[61,1] = 230
$data[61,2] = 0
$data[61,3] = 65
$data[61,4] = 148
$data[61,5] = 0
$data[61,6] = 0
$data[61,7] = 17
$data[62,0] = "Taiwan"
$data[62,1] = 216
$data[62,2] = 21
$data[62,3] = 29
$data[62,4] = 185
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 2
$data[63,0] = "Eslovaquia"
$data[63,1] = 204
$data[63,2] = 18
$data[63,3] = 7
$data[63,4] = 197
$data[63,5] = 2
$data[63,6] = 0
$data[63,7] = 0
$data[64,0] = "Bulgaria"
$data[64,1] = 202
$data[64,2] = 1
$data[64,3] = 3
$data[64,4] = 196
$data[64,5] = 8
$data[64,6] = 0
$data[64,7] = 3
$data[65,0] = "Emiratos Arabes Unidos"
$data[65,1] = 198
$data[65,2] = 0
$data[65,3] = 41
$data[65,4] = 155
$data[65,5] = 2
$data[65,6] = 0
$data[65,7] = 2
$data[66,0] = "Letonia"
$data[66,1] = 197
$data[66,2] = 17
$data[66,3] = 1
$data[66,4] = 196
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 0
$data[67,0] = "Kuwait"
$data[67,1] = 191
$data[67,2] = 2
$data[67,3] = 39
$data[67,4] = 152
$data[67,5] = 5
$data[67,6] = 0
$data[67,7] = 0
$data[68,0] = "Lituania"
$data[68,1] = 187
$data[68,2] = 8
$data[68,3] = 1
$data[68,4] = 185
$data[68,5] = 1
$data[68,6] = 0
$data[68,7] = 1
$data[69,0] = "San Marino"
$data[69,1] = 187
$data[69,2] = 0
$data[69,3] = 4
$data[69,4] = 162
$data[69,5] = 12
$data[69,6] = 1
$data[69,7] = 21
$data[70,0] = "Hungria"
$data[70,1] = 187
$data[70,2] = 20
$data[70,3] = 21
$data[70,4] = 157
$data[70,5] = 6
$data[70,6] = 1
$data[70,7] = 9
$data[71,0] = "Principado de Andorra"
$data[71,1] = 164
$data[71,2] = 31
$data[71,3] = 1
$data[71,4] = 162
$data[71,5] = 7
$data[71,6] = 0
$data[71,7] = 1
$data[72,0] = "Uruguay"
$data[72,1] = 162
$data[72,2] = 0
$data[72,3] = 0
$data[72,4] = 162
$data[72,5] = 3
$data[72,6] = 0
$data[72,7] = 0
$data[73,0] = "Costa Rica"
$data[73,1] = 158
$data[73,2] = 0
$data[73,3] = 2
$data[73,4] = 154
$data[73,5] = 2
$data[73,6] = 0
$data[73,7] = 2
$data[74,0] = "Nueva Zelanda"
$data[74,1] = 155
$data[74,2] = 53
$data[74,3] = 12
$data[74,4] = 143
$data[74,5] = 0
$data[74,6] = 0
$data[74,7] = 0
$data[75,0] = "Bosnia y Herzegovina"
$data[75,1] = 150
$data[75,2] = 14
$data[75,3] = 2
$data[75,4] = 146
$data[75,5] = 1
$data[75,6] = 1
$data[75,7] = 2
$data[76,0] = "Republica de Macedonia"
$data[76,1] = 148
$data[76,2] = 12
$data[76,3] = 1
$data[76,4] = 145
$data[76,5] = 1
$data[76,6] = 0
$data[76,7] = 2
$data[77,0] = "Marruecos"
$data[77,1] = 143
$data[77,2] = 0
$data[77,3] = 5
$data[77,4] = 134
$data[77,5] = 1
$data[77,6] = 0
$data[77,7] = 4
$data[78,0] = "Jordania"
$data[78,1] = 127
$data[78,2] = 0
$data[78,3] = 1
$data[78,4] = 126
$data[78,5] = 0
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = "Albania"
$data[79,1] = 123
$data[79,2] = 19
$data[79,3] = 10
$data[79,4] = 108
$data[79,5] = 2
$data[79,6] = 1
$data[79,7] = 5
$data[80,0] = "Vietnam"
$data[80,1] = 123
$data[80,2] = 0
$data[80,3] = 17
$data[80,4] = 106
$data[80,5] = 3
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = "Islas Feroe"
$data[81,1] = 122
$data[81,2] = 4
$data[81,3] = 23
$data[81,4] = 99
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 0
$data[82,0] = "Republica de Chipre"
$data[82,1] = 116
$data[82,2] = 0
$data[82,3] = 3
$data[82,4] = 112
$data[82,5] = 3
$data[82,6] = 0
$data[82,7] = 1
$data[83,0] = "Tunez"
$data[83,1] = 114
$data[83,2] = 25
$data[83,3] = 1
$data[83,4] = 110
$data[83,5] = 11
$data[83,6] = 0
$data[83,7] = 3
$data[84,0] = "Burkina Faso"
$data[84,1] = 114
$data[84,2] = 15
$data[84,3] = 7
$data[84,4] = 103
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 4
$data[85,0] = "Malta"
$data[85,1] = 110
$data[85,2] = 3
$data[85,3] = 2
$data[85,4] = 108
$data[85,5] = 1
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = "Moldavia"
$data[86,1] = 109
$data[86,2] = 0
$data[86,3] = 2
$data[86,4] = 106
$data[86,5] = 10
$data[86,6] = 0
$data[86,7] = 1
$data[87,0] = "Brunei"
$data[87,1] = 104
$data[87,2] = 13
$data[87,3] = 2
$data[87,4] = 102
$data[87,5] = 2
$data[87,6] = 0
$data[87,7] = 0
$data[88,0] = "Sri Lanka"
$data[88,1] = 101
$data[88,2] = 4
$data[88,3] = 2
$data[88,4] = 99
$data[88,5] = 2
$data[88,6] = 0
$data[88,7] = 0
$data[89,0] = "Ucrania"
$data[89,1] = 97
$data[89,2] = 24
$data[89,3] = 1
$data[89,4] = 93
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 3
$data[90,0] = "Camboya"
$data[90,1] = 87
$data[90,2] = 0
$data[90,3] = 2
$data[90,4] = 85
$data[90,5] = 1
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = "Azerbaiyan"
$data[91,1] = 87
$data[91,2] = 15
$data[91,3] = 10
$data[91,4] = 76
$data[91,5] = 6
$data[91,6] = 0
$data[91,7] = 1
$data[92,0] = "Senegal"
$data[92,1] = 86
$data[92,2] = 7
$data[92,3] = 8
$data[92,4] = 78
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = "Venezuela"
$data[93,1] = 84
$data[93,2] = 0
$data[93,3] = 15
$data[93,4] = 69
$data[93,5] = 2
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = "Oman"
$data[94,1] = 84
$data[94,2] = 18
$data[94,3] = 17
$data[94,4] = 67
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = "Bielorrusia"
$data[95,1] = 81
$data[95,2] = 0
$data[95,3] = 22
$data[95,4] = 59
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 0
$data[96,0] = "Reunion"
$data[96,1] = 75
$data[96,2] = 4
$data[96,3] = 1
$data[96,4] = 74
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = "Kazajistan"
$data[97,1] = 68
$data[97,2] = 6
$data[97,3] = 0
$data[97,4] = 68
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = "Georgia"
$data[98,1] = 67
$data[98,2] = 6
$data[98,3] = 9
$data[98,4] = 58
$data[98,5] = 1
$data[98,6] = 0
$data[98,7] = 0
$data[99,0] = "Camerun"
$data[99,1] = 66
$data[99,2] = 10
$data[99,3] = 2
$data[99,4] = 64
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = "Guadalupe"
$data[100,1] = 62
$data[100,2] = 0
$data[100,3] = 0
$data[100,4] = 61
$data[100,5] = 4
$data[100,6] = 0
$data[100,7] = 1
$data[101,0] = "Estado de Palestina"
$data[101,1] = 60
$data[101,2] = 1
$data[101,3] = 16
$data[101,4] = 44
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 0
$data[102,0] = "Martinica"
$data[102,1] = 53
$data[102,2] = 0
$data[102,3] = 0
$data[102,4] = 52
$data[102,5] = 7
$data[102,6] = 0
$data[102,7] = 1
$data[103,0] = "Trinidad yTobago"
$data[103,1] = 52
$data[103,2] = 1
$data[103,3] = 0
$data[103,4] = 52
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = "Ghana"
$data[104,1] = 52
$data[104,2] = 25
$data[104,3] = 0
$data[104,4] = 50
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 2
$data[105,0] = "Liechtenstein"
$data[105,1] = 51
$data[105,2] = 0
$data[105,3] = 0
$data[105,4] = 51
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 0
$data[106,0] = "Uzbekistan"
$data[106,1] = 50
$data[106,2] = 4
$data[106,3] = 0
$data[106,4] = 50
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 0
$data[107,0] = "Consejo Danes para los Refugiados"
$data[107,1] = 45
$data[107,2] = 9
$data[107,3] = 0
$data[107,4] = 43
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 2
$data[108,0] = "Kirguistan"
$data[108,1] = 42
$data[108,2] = 26
$data[108,3] = 0
$data[108,4] = 42
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 0
$data[109,0] = "Afganistan"
$data[109,1] = 42
$data[109,2] = 2
$data[109,3] = 1
$data[109,4] = 40
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 1
$data[110,0] = "Nigeria"
$data[110,1] = 42
$data[110,2] = 2
$data[110,3] = 2
$data[110,4] = 39
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 1
$data[111,0] = "Cuba"
$data[111,1] = 40
$data[111,2] = 0
$data[111,3] = 0
$data[111,4] = 39
$data[111,5] = 3
$data[111,6] = 0
$data[111,7] = 1
$data[112,0] = "Puerto Rico"
$data[112,1] = 39
$data[112,2] = 8
$data[112,3] = 1
$data[112,4] = 36
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 2
$data[113,0] = "Banglades"
$data[113,1] = 39
$data[113,2] = 6
$data[113,3] = 5
$data[113,4] = 30
$data[113,5] = 0
$data[113,6] = 1
$data[113,7] = 4
$data[114,0] = "Ruanda"
$data[114,1] = 36
$data[114,2] = 0
$data[114,3] = 0
$data[114,4] = 36
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 0
$data[115,0] = "Mayotte"
$data[115,1] = 36
$data[115,2] = 12
$data[115,3] = 0
$data[115,4] = 36
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = "Mauricio"
$data[116,1] = 36
$data[116,2] = 0
$data[116,3] = 0
$data[116,4] = 34
$data[116,5] = 1
$data[116,6] = 0
$data[116,7] = 2
$data[117,0] = "Guam"
$data[117,1] = 32
$data[117,2] = 3
$data[117,3] = 0
$data[117,4] = 31
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 1
$data[118,0] = "Honduras"
$data[118,1] = 30
$data[118,2] = 0
$data[118,3] = 0
$data[118,4] = 30
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = "Montenegro"
$data[119,1] = 29
$data[119,2] = 2
$data[119,3] = 0
$data[119,4] = 28
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 1
$data[120,0] = "Bolivia"
$data[120,1] = 28
$data[120,2] = 1
$data[120,3] = 0
$data[120,4] = 28
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = "Paraguay"
$data[121,1] = 27
$data[121,2] = 5
$data[121,3] = 0
$data[121,4] = 25
$data[121,5] = 1
$data[121,6] = 1
$data[121,7] = 2
$data[122,0] = "Kenia"
$data[122,1] = 25
$data[122,2] = 9
$data[122,3] = 0
$data[122,4] = 25
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = "Costa de Marfil"
$data[123,1] = 25
$data[123,2] = 0
$data[123,3] = 2
$data[123,4] = 23
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = "Macao"
$data[124,1] = 25
$data[124,2] = 0
$data[124,3] = 10
$data[124,4] = 15
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = "Polinesia Francesa"
$data[125,1] = 23
$data[125,2] = 5
$data[125,3] = 0
$data[125,4] = 23
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Monaco"
$data[126,1] = 23
$data[126,2] = 0
$data[126,3] = 1
$data[126,4] = 22
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = "Guayana Francesa"
$data[127,1] = 23
$data[127,2] = 3
$data[127,3] = 6
$data[127,4] = 17
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = "Jamaica"
$data[128,1] = 21
$data[128,2] = 2
$data[128,3] = 2
$data[128,4] = 18
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 1
$data[129,0] = "Isla de Man"
$data[129,1] = 20
$data[129,2] = 7
$data[129,3] = 0
$data[129,4] = 20
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = "Guatemala"
$data[130,1] = 20
$data[130,2] = 0
$data[130,3] = 0
$data[130,4] = 19
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 1
$data[131,0] = "Guyana"
$data[131,1] = 20
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 19
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 1
$data[132,0] = "Togo"
$data[132,1] = 20
$data[132,2] = 2
$data[132,3] = 1
$data[132,4] = 19
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = "Islas Virgenes de los Estados Unidos"
$data[133,1] = 17
$data[133,2] = 0
$data[133,3] = 0
$data[133,4] = 17
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Barbados"
$data[134,1] = 17
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 17
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = "Madagascar"
$data[135,1] = 17
$data[135,2] = 5
$data[135,3] = 0
$data[135,4] = 17
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = "Gibraltar"
$data[136,1] = 15
$data[136,2] = 0
$data[136,3] = 5
$data[136,4] = 10
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = "Maldivas"
$data[137,1] = 13
$data[137,2] = 0
$data[137,3] = 5
$data[137,4] = 8
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = "Etiopia"
$data[138,1] = 12
$data[138,2] = 1
$data[138,3] = 0
$data[138,4] = 12
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = "Tanzania"
$data[139,1] = 12
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 12
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = "Aruba"
$data[140,1] = 12
$data[140,2] = 3
$data[140,3] = 1
$data[140,4] = 11
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = "Mongolia"
$data[141,1] = 10
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 10
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = "Nueva Caledonia"
$data[142,1] = 10
$data[142,2] = 2
$data[142,3] = 0
$data[142,4] = 10
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = "Guinea Ecuatorial"
$data[143,1] = 9
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 9
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = "Uganda"
$data[144,1] = 9
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 9
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = "San Martin (Parte Francesa)"
$data[145,1] = 8
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 8
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = "Seychelles"
$data[146,1] = 7
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 7
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = "Haiti"
$data[147,1] = 6
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 6
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = "Surinam"
$data[148,1] = 6
$data[148,2] = 1
$data[148,3] = 0
$data[148,4] = 6
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = "Bermudas"
$data[149,1] = 6
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 6
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = "Benin"
$data[150,1] = 6
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 6
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = "Gabon"
$data[151,1] = 6
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 5
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 1
$data[152,0] = "El Salvador"
$data[152,1] = 5
$data[152,2] = 2
$data[152,3] = 0
$data[152,4] = 5
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = "Islas Caimanes"
$data[153,1] = 5
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 4
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 1
$data[154,0] = "Fiyi"
$data[154,1] = 4
$data[154,2] = 1
$data[154,3] = 0
$data[154,4] = 4
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = "Namibia"
$data[155,1] = 4
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 4
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Congo"
$data[156,1] = 4
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 4
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = "Groenlandia"
$data[157,1] = 4
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 4
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = "Bahamas"
$data[158,1] = 4
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 4
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = "Suazilandia"
$data[159,1] = 4
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 4
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = "Guinea"
$data[160,1] = 4
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 4
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = "Curazao"
$data[161,1] = 4
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 3
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 1
$data[162,0] = "Republica del Chad"
$data[162,1] = 3
$data[162,2] = 1
$data[162,3] = 0
$data[162,4] = 3
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = "Antigua y Barbuda"
$data[163,1] = 3
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 3
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = "Liberia"
$data[164,1] = 3
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 3
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Angola"
$data[165,1] = 3
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 3
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = "Niger"
$data[166,1] = 3
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 3
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Santa Lucia"
$data[167,1] = 3
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 3
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Zambia"
$data[168,1] = 3
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 3
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = "Republica de Africa Central"
$data[169,1] = 3
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 3
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Republica de Yibuti"
$data[170,1] = 3
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 3
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "San Bartolome"
$data[171,1] = 3
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 3
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "Sudan"
$data[172,1] = 3
$data[172,2] = 1
$data[172,3] = 0
$data[172,4] = 2
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 1
$data[173,0] = "Cabo Verde"
$data[173,1] = 3
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 2
$data[173,5] = 0
$data[173,6] = 1
$data[173,7] = 1
$data[174,0] = "Zimbabue"
$data[174,1] = 3
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 2
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 1
$data[175,0] = "San Martin (Parte Holandesa)"
$data[175,1] = 2
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 2
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "Butan"
$data[176,1] = 2
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 2
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = "Birmania"
$data[177,1] = 2
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 2
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Nicaragua"
$data[178,1] = 2
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 2
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = "Mauritania"
$data[179,1] = 2
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 2
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Dominica"
$data[180,1] = 2
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 2
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Laos"
$data[181,1] = 2
$data[181,2] = 2
$data[181,3] = 0
$data[181,4] = 2
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = "Nepal"
$data[182,1] = 2
$data[182,2] = 0
$data[182,3] = 1
$data[182,4] = 1
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = "Gambia"
$data[183,1] = 2
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 1
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 1
$data[184,0] = "Belice"
$data[184,1] = 1
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 1
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = "Somalia"
$data[185,1] = 1
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 1
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = "Papua Nueva Guinea"
$data[186,1] = 1
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 1
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = "Timor Oriental"
$data[187,1] = 1
$data[187,2] = 0
$data[187,3] = 0
$data[187,4] = 1
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = "Eritrea"
$data[188,1] = 1
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 1
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = "Mozambique"
$data[189,1] = 1
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 1
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = "Siria"
$data[190,1] = 1
$data[190,2] = 0
$data[190,3] = 0
$data[190,4] = 1
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = "San Vicente y las Granadinas"
$data[191,1] = 1
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 1
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = "Montserrat"
$data[192,1] = 1
$data[192,2] = 0
$data[192,3] = 0
$data[192,4] = 1
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Granada"
$data[193,1] = 1
$data[193,2] = 0
$data[193,3] = 0
$data[193,4] = 1
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = "Islas Turcas y Caicos"
$data[194,1] = 1
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 1
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = "Santa Sede"
$data[195,1] = 1
$data[195,2] = 0
$data[195,3] = 0
$data[195,4] = 1
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$ws.Range("A4:H199").Value = $data
